$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numeric-looking values stored as text (inlineStr) in the
# original workbook. To keep them as text (not auto-converted to numbers)
# we force a text number-format before assigning, then restore the style.
function Set-TextValue($cell, $value) {
    $cell.Style = "Normal"
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "268.04"
Set-TextValue $ws.Range("D3") "21.39"
Set-TextValue $ws.Range("D4") "6.196"
Set-TextValue $ws.Range("D5") "0.06165"
Set-TextValue $ws.Range("D6") "3.567"
Set-TextValue $ws.Range("D7") "6.519"
Set-TextValue $ws.Range("D8") "1.384"
Set-TextValue $ws.Range("D9") "0.8246"
Set-TextValue $ws.Range("D10") "0.01354"
Set-TextValue $ws.Range("D11") "0.1578"
Set-TextValue $ws.Range("D12") "0.08081"
Set-TextValue $ws.Range("D13") "0.03358"
Set-TextValue $ws.Range("D14") "0.03183"
Set-TextValue $ws.Range("D15") "0.09253"
Set-TextValue $ws.Range("D16") "3.767"
Set-TextValue $ws.Range("D18") "0.04674"
Set-TextValue $ws.Range("D19") "0.006327"
Set-TextValue $ws.Range("D20") "0.006206"
Set-TextValue $ws.Range("D24") "2.440"
Set-TextValue $ws.Range("D25") "0.3301"
Set-TextValue $ws.Range("D26") "0.1240"
Set-TextValue $ws.Range("D40") "0.04654"
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue $ws.Range("D42") "0.1128"
$ws.Range("E42").Value = "41BKEXTokenBKK"
$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue $ws.Range("D43") "0.003441"
$ws.Range("E43").Value = "42CEJICEJI"
Set-TextValue $ws.Range("D44") "0.01186"
Set-TextValue $ws.Range("D45") "0.00005826"
Set-TextValue $ws.Range("D46") "0.0009875"
Set-TextValue $ws.Range("D48") "0.7802"
Set-TextValue $ws.Range("D49") "0.002436"
Set-TextValue $ws.Range("D51") "0.01237"
